# Apply the "tests for tiaraoutputs and shielding output" edit to the
# "Experimental benchmarks" worksheet of the JADE mainconfig workbook.
#
# Rows affected (column D = MCNP, column H = Post-Processing):
#   Row 7  -> "Tiara Bonner Sphere detector"   (Tiara-BS)
#   Row 8  -> "Tiara Fission Chambers detector" (Tiara-FC)
#   Row 11 -> "FNG Tungsten"                    (FNG-W)
#
# For rows 7 and 8 both MCNP (D) and Post-Processing (H) are switched on
# ("false" -> "true"); for row 11 only MCNP (D) is switched on.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Experimental benchmarks")

# Cells that already contain the literal text "true" (e.g. D4/H4) use the
# worksheet's normal body style. Use one of them as the formatting template
# so the edited cells pick up the same (non quote-prefixed) style instead of
# the "quote prefix" style that some of the "false" cells use.
$trueTemplate = $ws.Range("D4")

function Set-TrueFlag($cellRef) {
    $cell = $ws.Range($cellRef)
    # Leading apostrophe forces this to be stored as literal text "true"
    # (matching the existing "true"/"false" shared strings) instead of
    # being auto-converted into a native Excel boolean.
    $cell.Value = "'true"
    $trueTemplate.Copy()
    $cell.PasteSpecial(-4122)  # xlPasteFormats
}

Set-TrueFlag "D7"
Set-TrueFlag "H7"
Set-TrueFlag "D8"
Set-TrueFlag "H8"
Set-TrueFlag "D11"

$excel.CutCopyMode = 0

# Update the selection shown on the "Experimental benchmarks" sheet.
$null = $ws.Activate()
$null = $ws.Range("H8").Select()
